$wb = $excel.ActiveWorkbook

$urls = @(
    "https://www.amazon.com/Holy-Stone-Quadcopter-Adjustable-Intelligent/dp/B074YYVXQH/ref=sr_1_4?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-4&keywords=drone",
    "https://www.amazon.com/Camera-EACHINE-Quadcopter-Wide-angle-Foldable/dp/B0776QJNS3/ref=sr_1_5?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-5&keywords=drone",
    "https://www.amazon.com/Holy-Stone-Predator-Helicopter-Quadcopter/dp/B0157IHJMQ/ref=sr_1_6?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-6&keywords=drone",
    "https://www.amazon.com/Holy-Stone-Shadow-Quadcopter-Beginners/dp/B074S2HK59/ref=sr_1_7?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-7&keywords=drone",
    "https://www.amazon.com/DROCON-Beginners-Training-Quadcopter-Operation/dp/B073HYDPT3/ref=sr_1_8?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-8&keywords=drone",
    "https://www.amazon.com/Cheerwing-X5SW-V3-Explorers2-Headless-Quadcopter/dp/B011JV9HA2/ref=sr_1_10?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-10&keywords=drone",
    "https://www.amazon.com/LBLA-Headless-Quadcopter-Compatible-Headset/dp/B077BZQ8JZ/ref=sr_1_12?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-12&keywords=drone",
    "https://www.amazon.com/Holy-Stone-Quadcopter-Altitude-Function/dp/B00SAUAP5C/ref=sr_1_16?s=toys-and-games&ie=UTF8&qid=1527886567&sr=1-16&keywords=drone"
)

$rows = @(2, 3, 4, 5, 6, 8, 10, 11)

foreach ($sheetName in @("PC", "drone")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $rows[$i]
        $ws.Range("B$r").Value = $urls[$i]
    }
}
